$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

$ws1.Range("F2").Value = 1520
$ws1.Range("F4").Value = 2113
$ws1.Range("F5").Value = 7825
$ws1.Range("F6").Value = 4883
$ws1.Range("F7").Value = 7174
$ws1.Range("F8").Value = 295
$ws1.Range("F9").Value = 1527
$ws1.Range("F11").Value = 212
$ws1.Range("F12").Value = 71
$ws1.Range("F13").Value = 1187
$ws1.Range("F16").Value = 24
$ws1.Range("F20").Value = 1234
$ws1.Range("F22").Value = 556
$ws1.Range("F24").Value = 1262
$ws1.Range("F25").Value = 52
$ws1.Range("F26").Value = 163
$ws1.Range("F28").Value = 17
$ws1.Range("F30").Value = 221
$ws1.Range("F31").Value = 897
$ws1.Range("F34").Value = 153
$ws1.Range("F35").Value = 134
$ws1.Range("F37").Value = 557
$ws1.Range("F38").Value = 571
$ws1.Range("F40").Value = 90
$ws1.Range("F43").Value = 428
$ws1.Range("F45").Value = 601
$ws1.Range("F46").Value = 160
$ws2.Range("F22").Value = 145
$ws2.Range("F29").Value = 32
$ws2.Range("F38").Value = 121
$ws2.Range("F44").Value = 15
$ws3.Range("F6").Value = 702
$ws3.Range("F7").Value = 184
$ws3.Range("F9").Value = 1738
$ws3.Range("F10").Value = 2652
$ws4.Range("F3").Value = 1520
$ws4.Range("F6").Value = 702
$ws4.Range("F7").Value = 7825
$ws4.Range("F8").Value = 184
$ws4.Range("F9").Value = 4883
$ws4.Range("F10").Value = 7174
$ws4.Range("F11").Value = 295
$ws4.Range("F12").Value = 1527
$ws4.Range("F15").Value = 212
$ws4.Range("F16").Value = 1738
$ws4.Range("F17").Value = 2652
$ws4.Range("F19").Value = 1187
$ws4.Range("F21").Value = 24
$ws4.Range("F23").Value = 1234
$ws4.Range("F26").Value = 1263
$ws4.Range("F27").Value = 163
$ws4.Range("F28").Value = 17
$ws4.Range("F29").Value = 221
$ws4.Range("F32").Value = 32
$ws4.Range("F34").Value = 153
$ws4.Range("F36").Value = 134
$ws4.Range("F38").Value = 571
$ws4.Range("F40").Value = 90
$ws4.Range("F42").Value = 121
$ws4.Range("F43").Value = 428
$ws4.Range("F44").Value = 601
$ws4.Range("F46").Value = 160
